# Rename the two Pearson-logo pictures (in the section footers) from
# "image1.png" to "image2.png", and the two BTec-logo pictures
# (in the section headers) from "image2.jpg" to "image1.jpg".
#
# The pictures are distinguished by their (unchanged) AlternativeText /
# description, which survives a round trip, rather than by a hard-coded
# Headers/Footers index - that keeps this correct regardless of which
# physical header/footer part ("default" vs "first page") Word exposes
# as Item(1) vs Item(2).

$d = $word.ActiveDocument

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    for ($i = 1; $i -le 3; $i++) {

        $h = $sec.Headers.Item($i)
        if ($h.Exists) {
            $shapes = $h.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $shp = $shapes.Item($j)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image1.jpg"
                }
            }
        }

        $f = $sec.Footers.Item($i)
        if ($f.Exists) {
            $shapes = $f.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $shp = $shapes.Item($j)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image2.png"
                }
            }
        }
    }
}
